$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Odabir automobila")

# Fill the "utility" formulas in rows 34-38 (columns D-I), mirroring the
# raw criteria values in rows 21-25, transformed by the max/min rule
# declared in row 19.
$srcRow = 21
for ($r = 34; $r -le 38; $r++) {
    foreach ($col in @("D","E","F","G","H","I")) {
        $formula = '=IF(' + $col + '$19="max",' + $col + $srcRow + ',1/' + $col + $srcRow + ')'
        $ws.Range($col + $r).Formula = $formula
    }
    $srcRow++
}

# Row 38 previously had empty cells with a bottom border (style index 2);
# once filled in, it matches the borderless style used by the rows above
# (style index 1), so the bottom border needs to be removed.
$ws.Range("D38:I38").Borders.Item(9).LineStyle = -4142

# Update the active selection left by the author on this sheet.
$ws.Range("Q32").Select()
